$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels (upper-cased entity column names)
$ws.Range("A2").Value = "ID"
$ws.Range("B2").Value = "MODEL_Name"
$ws.Range("C2").Value = "MODEL_ART"
$ws.Range("D2").Value = "Model_tYp"

# Column H briefly mirrors column C's text so AutoFit leaves behind the
# same stale <col bestFit> metadata seen in the target file, then is
# cleared along with the other now-unused cells.
$ws.Range("H2").Value = "MODEL_ART"
$ws.Range("E2").Value = "modEL_NuMBER"

$ws.Columns("B:B").AutoFit() | Out-Null
$ws.Columns("C:C").AutoFit() | Out-Null
$ws.Columns("E:E").AutoFit() | Out-Null
$ws.Columns("H:H").AutoFit() | Out-Null

$ws.Range("F2").Value = $null
$ws.Range("H2").Value = $null
$ws.Range("I2").Value = $null
$ws.Range("L2").Value = $null

$ws.Range("E5").Select() | Out-Null
